$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("H9").Value = 3
$ws.Range("D21").Value = 11.01333333333333
$ws.Range("F21").Value = 11.13
$ws.Range("D25").Value = 7.413333333333334
$ws.Range("F25").Value = 10.43
$ws.Range("D29").Value = 11.83833333333333
$ws.Range("F29").Value = 11.655
$ws.Range("D33").Value = 7.988333333333333
$ws.Range("F33").Value = 10.605
$ws.Range("H40").Value = 2
$ws.Range("D41").Value = 14.28833333333333
$ws.Range("F41").Value = 14.205
$ws.Range("D43").Value = 15.54833333333333
$ws.Range("F43").Value = 14.765
$ws.Range("D49").Value = 9.763333333333334
$ws.Range("F49").Value = 10.78
$ws.Range("H49").Value = 1
$ws.Range("D50").Value = 13.30333333333333
$ws.Range("F50").Value = 14.52
$ws.Range("D52").Value = 14.03833333333333
$ws.Range("F52").Value = 10.255
$ws.Range("D59").Value = 15.54833333333333
$ws.Range("F59").Value = 13.365
$ws.Range("D60").Value = 16.96833333333333
$ws.Range("F60").Value = 16.685
$ws.Range("D68").Value = 12.36166666666667
$ws.Range("F68").Value = 11.27
$ws.Range("D81").Value = 10.78166666666667
$ws.Range("F81").Value = 10.29
$ws.Range("D86").Value = 8.323333333333331
$ws.Range("F86").Value = 12.74
$ws.Range("D94").Value = 17.75833333333333
$ws.Range("F94").Value = 19.975
$ws.Range("D100").Value = 10.96166666666667
$ws.Range("F100").Value = 9.52
$ws.Range("D107").Value = 10.09833333333333
$ws.Range("F107").Value = 10.22
$ws.Range("D108").Value = 6.563333333333333
$ws.Range("F108").Value = 7.035
$ws.Range("D119").Value = 8.561666666666666
$ws.Range("F119").Value = 9.17
$ws.Range("D124").Value = 8.133333333333333
$ws.Range("F124").Value = 5.949999999999999
$ws.Range("D125").Value = 0
$ws.Range("F125").Value = 0
$ws.Range("H127").Value = 7
$ws.Range("D141").Value = 1.208333333333331
$ws.Range("F141").Value = 5.425
$ws.Range("D147").Value = 4.413333333333334
$ws.Range("F147").Value = 3.43
$ws.Range("H149").Value = 6
$ws.Range("D172").Value = 8.283333333333331
$ws.Range("F172").Value = 7.899999999999999
$ws.Range("D193").Value = 1.303333333333331
$ws.Range("F193").Value = 6.02
$ws.Range("D210").Value = 5.153333333333332
$ws.Range("F210").Value = 4.864999999999999
$ws.Range("D226").Value = 5.208333333333333
$ws.Range("F226").Value = 5.109999999999999
$ws.Range("D241").Value = 0
$ws.Range("F241").Value = 0
$ws.Range("D248").Value = 2.078333333333333
$ws.Range("F248").Value = 1.61
$ws.Range("D250").Value = 5.458333333333331
$ws.Range("F250").Value = 6.675
$ws.Range("D251").Value = 5.323333333333331
$ws.Range("F251").Value = 11.54
$ws.Range("D256").Value = 14.04333333333333
$ws.Range("F256").Value = 11.16
$ws.Range("D257").Value = 0
$ws.Range("F257").Value = 0
$ws.Range("D265").Value = 6.538333333333331
$ws.Range("F265").Value = 3.255
$ws.Range("D285").Value = 0
$ws.Range("F285").Value = 10
$ws.Range("D293").Value = 5.528333333333333
$ws.Range("F293").Value = 5.345
$ws.Range("D307").Value = 7.306666666666665
$ws.Range("F307").Value = 13.52333333333333
$ws.Range("D315").Value = 0.4833333333333314
$ws.Range("F315").Value = 0.7
$ws.Range("D320").Value = 0
$ws.Range("F320").Value = 0
$ws.Range("D328").Value = 1.453333333333334
$ws.Range("F328").Value = 4.47
$ws.Range("D341").Value = 9.924166666666665
$ws.Range("F341").Value = 9.532499999999999
$ws.Range("D361").Value = -2.816666666666669
$ws.Range("F361").Value = 1.4
$ws.Range("D428").Value = 5.616666666666665
$ws.Range("F428").Value = 6.225
$ws.Range("D434").Value = 0.4716666666666667
$ws.Range("F434").Value = 0.3733333333333333
$ws.Range("D468").Value = 6.924999999999999
$ws.Range("F468").Value = 5.033333333333333
$ws.Range("H468").Value = 1
$ws.Range("D496").Value = 0
$ws.Range("D525").Value = -4.246666666666666
$ws.Range("F525").Value = 0.06999999999999999
$ws.Range("D543").Value = 0.9233333333333321
$ws.Range("F543").Value = 0.84
$ws.Range("D585").Value = 0
$ws.Range("H588").Value = 4
$ws.Range("H597").Value = 7
$ws.Range("H598").Value = 4
$ws.Range("H651").Value = 7
$ws.Range("H660").Value = 7
$ws.Range("D663").Value = 0
$ws.Range("H677").Value = 6
$ws.Range("D678").Value = -4.616666666666669
$ws.Range("F678").Value = -1.4

Write-Host "Applied 113 cell updates"
